$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(0.1630846444411361, -1.088872551767469, 0.2388572144868575, -0.1694399950501343, 0.5040425926957953, -0.1932960119371557, 0.04367710249015017, 0.07925193655196383, 0.9052087946043501, 0.2491603105793238),
    @(0.6966127121850244, 0.03971143396413679, 0.6050679262073962, -0.1390344077924999, 0.0777199357846518, 0.1045584651549948, 0.9267429767877318, 0.2690667718748542, -0.3899018488887417, 0.008856469899756692),
    @(0.5889366529896609, -0.1629753181314043, 0.05339438235456767, 0.08260390735964929, 0.906862924893769, 0.2504397716903112, -0.4078772338215886, -0.00880507520269791, -0.0702522111134038, -0.5044206789217187),
    @(0.135347353279659, 0.0918012475915978, 0.8871553623927978, 0.2212014561932883, -0.4400216186372514, -0.04180104165272297, -0.1034922560287353, -0.5377297743012852, 0.4410516555687047, 0.2139302114848279),
    @(0.8583540716299061, 0.219077672514715, -0.4327386538548288, -0.03138452870947878, -0.09205111059441351, -0.5259549742272362, 0.452935276392097, 0.2258494563263944, -0.4076000933034996, 0.0198003091970762),
    @(-0.5546208149506421, -0.1232296976428027, -0.1719238571658218, -0.6008843017759042, 0.3800526390729398, 0.1538132908583271, -0.4792866254931878, -0.05174196790978589, -0.2180182939146202, -0.09528239010531198),
    @(-0.07660908550845907, -0.5226165970886194, 0.4506431296588137, 0.2208746131534756, -0.4138532542967197, 0.01293740259710929, -0.153689426264273, -0.03111699931576278, 0.01605047729302533, -0.3805999332580445),
    @(0.7107656714615536, 0.3617805447359003, -0.3284087146615119, 0.07258913989446059, -0.106032695953167, 0.01096152432303771, 0.0555349373482778, -0.3423217815153134, -0.0888458464935532, 0.1310104539444903),
    @(-0.5260169647433546, -0.05137068376835741, -0.1974123165853111, -0.06599669411565129, -0.01502993377330319, -0.410047744940668, -0.1553091497739889, 0.06510967366033998, -0.1705525061135018, -0.3494711127835084),
    @(-0.1300102632924844, 0.00888099880676374, 0.06313616294635133, -0.3304359561846383, -0.07506199009359621, 0.1456359881641034, -0.08990357932170889, -0.2687683468051136, -0.2792820200016274, -0.4869807987902649),
    @(0.05077946203237782, -0.344626946110449, -0.09000824431623794, 0.1303795298704019, -0.1052864881356162, -0.2842023355090271, -0.2947364119262331, -0.5024432253954677, -0.2852286099811868, 0.02446021624927014),
    @(0.0822944582731665, 0.2248716626294139, -0.04696104924847117, -0.2426806063420137, -0.2610220372837692, -0.4723563104746108, -0.2568270884946112, 0.05207866879825723, -0.1025880601374284, 0.5476698920998867),
    @(-0.1414695961309799, -0.2791950053566942, -0.2708243371007503, -0.4698642071905083, -0.2486761060184952, 0.06283442082819635, -0.09063329381595758, 0.560176603242597, 0.3441873680460773, -0.06276115523468917),
    @(-0.1284196148407921, -0.3892880536239818, -0.1956617690156243, 0.1035623046680138, -0.05538318383206281, 0.592984168753176, 0.3759056446645694, -0.03152873430779118, 0.4673336997040285, 0.3084050283952137),
    @(-0.01772937608965819, 0.1904199002428446, -0.01245745804686382, 0.6148928434241717, 0.3877598625592121, -0.02448210058422801, 0.4720825582591874, 0.3120561101207888, 0.3676602717569419, 2.617197116239002),
    @(0.2009879203755026, -0.003682669475092215, 0.622695247467907, 0.3950699233937206, -0.01742125293059998, 0.4790177347602291, 0.3189281133539701, 0.3745006084643348, 2.624021620117335, 10.09303763079093),
    @(-0.1106570732907728, 0.5574356245296381, 0.3494793248607413, -0.0537519294468537, 0.4470476400263193, 0.2890118858004561, 0.3455519702652192, 2.59552891386994, 10.06475980582587, -8.106464717219385),
    @(0.5913088222741343, 0.3522324851175385, -0.06534086029760827, 0.4289262259694052, 0.267911772437701, 0.3230936335787457, 2.57245120070562, 10.04139962167697, -8.129953734994666, 0.0096762954492835),
    @(0.09495705561128835, -0.1966201317798049, 0.3555705446956575, 0.2206704977705274, 0.2876467492000076, 2.542330314111783, 10.01368383765935, -8.156583427754628, -0.01646294162652273, 2.131310184089503),
    @(-0.2398842457485278, 0.3272291486687826, 0.1959944722584051, 0.2633601162231758, 2.517707017856339, 9.988700664718653, -8.181813387285533, -0.04183851343822592, 2.105855004056436, -1.334125068513618),
    @(0.4374203038949512, 0.2521171062027519, 0.2974419321570869, 2.542806401728791, 10.0101393154344, -8.161866677895206, -0.02249986200433229, 2.124945828975858, -1.315135252606165, -1.376330353346725),
    @(0.09526702093765275, 0.2022742013914471, 2.48837927106286, 9.97358532802263, -8.18953186637841, -0.04590700167627143, 2.10360890909983, -1.335467850351198, -1.396174332378672, 0.6455109801188312),
    @(0.1534105877790901, 2.460610420786189, 9.956365482985138, -8.201626982332003, -0.05549112916037019, 2.095251978282879, -1.343224567749304, -1.403637555454572, 0.6381912792035592, 0.05577943355819903)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}
